# Apply the edits captured in the commit: updated threshold values on
# Sheet1 (C3, C4, C5) plus the new selection left behind in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates -------------------------------------------------
$ws.Range("C3").Value = 9.5
$ws.Range("C4").Value = 1.4
$ws.Range("C5").Value = 25

# --- Selection update -----------------------------------------------
# The saved file shows the user had selected B2:C5 (with C5 as the last
# cell touched). Select the full range so the persisted sqref matches.
$ws.Range("B2:C5").Select()
